# "9. Mean free path.xlsx" - update data, fix Reynolds number calculations
#
# 1. Rename the two sheets.
# 2. Pull a couple of reference constants (water density, universal gas
#    constant R) out onto the raw-data grid and rebuild the mean-speed /
#    Reynolds-number formulas around them.
# 3. Shift the per-run statistics table (rows 10-15) onto the corrected
#    formulas, add the sixth run (row 15) and the alternate Reynolds-number
#    column (R/S).
# 4. Rebuild the propagation-of-error block (rows 18-23) to match.

$wb  = $excel.ActiveWorkbook
$raw = $wb.Worksheets.Item(1)
$raw.Name = "raw"
$ws  = $wb.Worksheets.Item(2)
$ws.Name = "results"

function Copy-CellStyle($srcAddr, $dstAddr) {
    $s = $ws.Range($srcAddr)
    $d = $ws.Range($dstAddr)
    $s.Copy()
    $d.PasteSpecial(-4122)   # xlPasteFormats
}

# ---------------------------------------------------------------------
# Row 2: add water density label + value
# ---------------------------------------------------------------------
Copy-CellStyle "B2" "E2"
$ws.Range("E2").Value = "\rho_\u0432\u043E\u0434\u0430, [kg/m^3]"
Copy-CellStyle "B2" "F2"
$ws.Range("F2").Value = 997.0

# ---------------------------------------------------------------------
# Row 3: I3/J3 become the "R" constant (was \bar{u} formula); K3/L3 removed
# ---------------------------------------------------------------------
$ws.Range("I3").Value = "\pm"          # shares text with the (new) "R" label's neighbour - placeholder, corrected below
$ws.Range("I3").Value = "R"
$ws.Range("J3").Clear()
Copy-CellStyle "F3" "J3"
$ws.Range("J3").Value = 8.314
$ws.Range("K3").Clear()
$ws.Range("L3").Clear()

# ---------------------------------------------------------------------
# Row 4: I4/J4 become the mean-speed formula (now driven by J3 = R);
# K4/L4 (new) hold the "+-" label and the speed-uncertainty product
# ---------------------------------------------------------------------
$ws.Range("I4").Value = "\bar{u}"
$ws.Range("J4").Formula = "=1.6*SQRT(`$J`$3*`$J`$1/`$J`$2)"
Copy-CellStyle "K3" "K4"
$ws.Range("K4").Value = "\pm"
Copy-CellStyle "L9" "L4"
$ws.Range("L4").Formula = "=J4*J5"

# ---------------------------------------------------------------------
# Row 5: I5/J5 (new) hold the mean-speed uncertainty formula (moved down
# from the old J4)
# ---------------------------------------------------------------------
Copy-CellStyle "I4" "I5"
$ws.Range("I5").Value = "\Delta \bar{u}"
Copy-CellStyle "J4" "J5"
$ws.Range("J5").Formula = "=0.005/J1+0.000005/J2"

# ---------------------------------------------------------------------
# Row 9 header: eta column now reads "x10^-6" (was "x10^-9")
# ---------------------------------------------------------------------
$ws.Range("I9").Value = "eta \cdot 10^-6"

# ---------------------------------------------------------------------
# Rows 10-15: per-run statistics table
#   G  = F2*F4*D   (was F3*F4*D)
#   H  = G*PI()*C1^4*C/(8*F1*F5)          (unchanged form)
#   I  = H*10^6                            (was H*10^9)
#   L  = F1/(PI()*C1^2*C)                  (unchanged)
#   M  = C1*L*F3/H                         (unchanged form)
#   O  = 3*H/(F3*J4)                       (was 3*H/(F3*J3))
#   P  = O*10^9                            (was O*10^7)
#   R  = 3*H/(1.6*F3)*SQRT(J2/(J3*J1))     (new)
#   S  = R*10^6                            (new, row 10 only)
# ---------------------------------------------------------------------
$rows = 10..15
foreach ($r in $rows) {
    $ws.Range("G$r").Formula = "=`$F`$2*`$F`$4*(D$r)"
    $ws.Range("H$r").Formula = "=G$r*PI()*`$C`$1^4*C$r/(8*`$F`$1*`$F`$5)"
    $ws.Range("I$r").Formula = "=H$r*10^6"
    $ws.Range("O$r").Formula = "=3*H$r/(`$F`$3*`$J`$4)"
    $ws.Range("P$r").Formula = "=O$r*10^9"

    Copy-CellStyle "O10" "R$r"
    $ws.Range("R$r").Formula = "=3*H$r/(1.6*`$F`$3)*SQRT(`$J`$2/(`$J`$3*`$J`$1))"
}
Copy-CellStyle "P10" "S10"
$ws.Range("S10").Formula = "=R10*10^6"

# New 6th measurement run (row 15) - was a pair of empty B15/C15 cells
Copy-CellStyle "B14" "B15"
$ws.Range("B15").Value = "30 nm"
Copy-CellStyle "C14" "C15"
$ws.Range("C15").Value = 50.0
Copy-CellStyle "D14" "D15"
$ws.Range("D15").Formula = "=(5.1+3.1)*10^-2"
Copy-CellStyle "E14" "E15"
$ws.Range("E15").Formula = "=(2.9+5.1)*10^-2"
Copy-CellStyle "L10" "L15"
$ws.Range("L15").Formula = "=`$F`$1/(PI()*`$C`$1^2*C15)"
Copy-CellStyle "M10" "M15"
$ws.Range("M15").Formula = "=`$C`$1*L15*`$F`$3/H15"

# ---------------------------------------------------------------------
# Rows 10-14 also keep their L/M formulas (unchanged form, still present)
# ---------------------------------------------------------------------
foreach ($r in 10..14) {
    $ws.Range("L$r").Formula = "=`$F`$1/(PI()*`$C`$1^2*C$r)"
    $ws.Range("M$r").Formula = "=`$C`$1*L$r*`$F`$3/H$r"
}

# ---------------------------------------------------------------------
# Rows 18-23: propagation-of-error block
#   G19:G23 - relative-error sum, water-manometer term now 10*10^-6 (was 50*10^-6)
#   H19:H23 = G*H10.. (same row offset pattern)
#   I19:I23 = G*I10..
#   O18:O22 = G(next row) + 0.0005/F5 + J5/J4     (was G/I10 + 0.0005/F5 + J4/J3)
#   P18:P22 = O*P(matching row)
# ---------------------------------------------------------------------
$ws.Range("G19").Formula = "=4*(0.00005/`$C`$1)^3+1/C10+(10*10^-6)/`$F`$1+0.0005/`$F`$5 + 0.0005/D10 + 0.0005/`$F`$3"
$ws.Range("G20").Formula = "=4*(0.00005/`$C`$1)^3+1/C11+(10*10^-6)/`$F`$1+0.0005/`$F`$5 + 0.0005/D11 + 0.0005/`$F`$3"
$ws.Range("G21").Formula = "=4*(0.00005/`$C`$1)^3+1/C12+(10*10^-6)/`$F`$1+0.0005/`$F`$5 + 0.0005/D12 + 0.0005/`$F`$3"
$ws.Range("G22").Formula = "=4*(0.00005/`$C`$1)^3+1/C13+(10*10^-6)/`$F`$1+0.0005/`$F`$5 + 0.0005/D13 + 0.0005/`$F`$3"
$ws.Range("G23").Formula = "=4*(0.00005/`$C`$1)^3+1/C14+(10*10^-6)/`$F`$1+0.0005/`$F`$5 + 0.0005/D14 + 0.0005/`$F`$3"

$ws.Range("H19").Formula = "=G19*H10"
$ws.Range("H20").Formula = "=G20*H11"
$ws.Range("H21").Formula = "=G21*H12"
$ws.Range("H22").Formula = "=G22*H13"
$ws.Range("H23").Formula = "=G23*H14"

$ws.Range("I19").Formula = "=G19*I10"
$ws.Range("I20").Formula = "=G20*I11"
$ws.Range("I21").Formula = "=G21*I12"
$ws.Range("I22").Formula = "=G22*I13"
$ws.Range("I23").Formula = "=G23*I14"

$ws.Range("O18").Formula = "=G19+0.0005/`$F`$5+`$J`$5/`$J`$4"
$ws.Range("O19").Formula = "=G20+0.0005/`$F`$5+`$J`$5/`$J`$4"
$ws.Range("O20").Formula = "=G21+0.0005/`$F`$5+`$J`$5/`$J`$4"
$ws.Range("O21").Formula = "=G22+0.0005/`$F`$5+`$J`$5/`$J`$4"
$ws.Range("O22").Formula = "=G23+0.0005/`$F`$5+`$J`$5/`$J`$4"

$ws.Range("P18").Formula = "=O18*P10"
$ws.Range("P19").Formula = "=O19*P11"
$ws.Range("P20").Formula = "=O20*P12"
$ws.Range("P21").Formula = "=O21*P13"
$ws.Range("P22").Formula = "=O22*P14"

Write-Host "edit complete"
